# Applies the Jan 14 2024 cryptos-list refresh (prices, 1h % changes,
# and a re-ranking that swaps TRON/Polkadot and replaces Algorand with Aave).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h change) hold numeric-looking text
# (e.g. "42.811.76", "  +0.21%  ") that must stay plain text, exactly as
# scraped -- force the Text number format before writing each value so
# Excel does not silently coerce them into numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.811.76'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.539.55'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.57'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.56'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +5.58%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.577'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.543'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.29'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0824'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.23%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.114'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.56'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.925.82'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.564.50'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.02'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +5.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.870'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.824.73'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.19'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0986'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.56'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.59'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '253.90'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.07'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.81'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.15'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.61%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +7.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.87'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.14'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '157.12'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.39'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +13.32%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.31'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0794'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.63'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -4.65%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '25.11'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.11%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.15'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +31.40%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.096.38'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0304'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.23%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.96'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.92'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.784.71'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.69'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +6.64%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '102.95'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.79%  '
